$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.265.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.365.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.21%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.365.65"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.06%  "

# Row 9
$ws.Range("E9").Value = "  -3.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.55%  "

# Row 11
$ws.Range("E11").Value = "  -4.30%  "

# Row 12
$ws.Range("E12").Value = "  -3.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.946.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.90%  "

# Row 14
$ws.Range("E14").Value = "  -1.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.370.18"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "  -6.18%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.334.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.36%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.22%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.11%  "

# Row 23
$ws.Range("E23").Value = "  -4.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.500.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.05%  "

# Row 25
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.63%  "

# Row 27
$ws.Range("E27").Value = "  -11.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.84%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.05%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.86%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.65%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.395.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "

# Row 35
$ws.Range("E35").Value = "  -6.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.95%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.49%  "

# Row 40
$ws.Range("E40").Value = "  -5.69%  "

# Row 41
$ws.Range("E41").Value = "  -5.00%  "

# Row 42
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.38%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.761"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.13%  "

# Row 46
$ws.Range("E46").Value = "  -7.23%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.87%  "

# Row 48
$ws.Range("E48").Value = "  -8.85%  "

# Row 49
$ws.Range("E49").Value = "  -3.45%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.233.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.04%  "
